$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 4 (the row currently holding account 005064129 / THIAGO)
# so the new record (account 005002457 / ROSANGELA) is added to the table,
# pushing THIAGO and all following rows down by one.
$ws.Rows.Item(4).Insert()

# Account numbers are zero-padded strings, not numbers, so force the
# column to Text before assigning - otherwise Excel would coerce
# "005002457" into the number 5002457 and the leading zeros would be lost.
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "005002457"

$ws.Cells.Item(4, 2).Value = "ROSANGELA"
$ws.Cells.Item(4, 3).Value = 33043.39
